# "Debugged + made some changes"
# Slide 6: the content placeholder's paragraph about training a classifier
# used to read "...train an instance of AdaBoost (Gradient Boosting)".
# Replace it with "...train an SVM".

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(6)
$shape = $slide.Shapes.Item(1)
$textRange = $shape.TextFrame.TextRange

$paragraph = $textRange.Paragraphs(5)
$paragraph.Text = "The features are then used to train an SVM"
